# This sheet is a weekly/daily price log for "Pomelo" (grapefruit) at the
# Vega Central Mapocho de Santiago market. A new daily record is being
# inserted as row 35 (dated serial 44803), pushing all subsequent rows
# (old 35..83) down by one (new 36..84). The dimension grows from
# A1:T83 to A1:T84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; existing rows 35:83 shift down to 36:84.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new price record. The
# non-numeric/reference columns (A, B, C, E, F, G, H, I, J, L, R, T) reuse
# the same values as the surrounding "Start Ruby" / "Primera" / Región
# Metropolitana records, while the date, variety, volume, prices, unit and
# $/Kg columns hold the new data point.
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44803
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100102
$ws.Range("H35").Value = "Cítricos"
$ws.Range("I35").Value = 100102006
$ws.Range("J35").Value = "Pomelo"
$ws.Range("K35").Value = "Start Ruby"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 280
$ws.Range("N35").Value = 12000
$ws.Range("O35").Value = 12000
$ws.Range("P35").Value = 12000
$ws.Range("Q35").Value = "$/caja 14 kilos"
$ws.Range("R35").Value = "Región Metropolitana"
$ws.Range("S35").Value = 857
$ws.Range("T35").Value = 14
